$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This script performs three related edits to the "Book of Gods" review doc:
#
# 1. Insert a new paragraph right after the title (Heading1) paragraph that
#    contains a bold "Meta description" label followed by the (previously
#    existing) meta-description sentence.
# 2. Delete the now-duplicated bold paragraph "Play Book of Gods for Free -
#    Review by Slot Expert" that used to sit near the end of the document
#    (right before the italic meta-description paragraph).
# 3. Replace the text of the remaining italic paragraph (which used to hold
#    the meta description) with a new AI-image-generation prompt describing
#    a Maya warrior illustration.
# ---------------------------------------------------------------------------

$headingText = "Play Book of Gods for Free - Review by Slot Expert"
$metaSentence = "Read our review of Book of Gods - an Ancient Egyptian-themed online slot game. Play this visually amazing slot for free and explore its exciting features."
$newImagePrompt = 'Create a cartoon-style image of a happy Maya warrior with glasses fitting the game "Book of Gods". The warrior should be holding the magic book with the Eye of Horus in the center and standing in front of a pyramid. The background should have a bright and vibrant color scheme, with symbols from Ancient Egypt surrounding the pyramid. The warrior should have a big smile on their face, showcasing their excitement at discovering the treasure hidden within the book. Make sure the image is eye-catching and engaging, depicting the sense of adventure and mystery that the game offers its players.'

# --- Locate, by content, the duplicated heading paragraph and the italic
#     meta-description paragraph near the end of the document. Paragraph
#     text always carries a trailing paragraph-mark character, so trim it
#     before comparing. ---
$total = $d.Paragraphs.Count
$dupIndex = -1
$italicIndex = -1
for ($i = 1; $i -le $total; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd()
    if ($t -eq $headingText -and $i -ne 1) {
        $dupIndex = $i
    }
    if ($t -eq $metaSentence) {
        $italicIndex = $i
    }
}

# --- Step 1: update the italic paragraph's text in place (keeps its
#     existing run/paragraph formatting, e.g. the <w:i/> run property). ---
if ($italicIndex -ne -1) {
    $italicPara = $d.Paragraphs.Item($italicIndex)
    $italicRange = $d.Range($italicPara.Range.Start, $italicPara.Range.End)
    $italicRange.Text = $newImagePrompt
}

# --- Step 2: delete the duplicated bold heading paragraph entirely
#     (including its paragraph mark, so no empty paragraph is left behind). ---
if ($dupIndex -ne -1) {
    $d.Paragraphs.Item($dupIndex).Range.Delete()
}

# --- Step 3: insert the new "Meta description" paragraph right after the
#     first paragraph (the document title / Heading1). ---
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"
$metaFullText = "Meta description: " + $metaSentence
$metaPara.Range.Text = $metaFullText

# Bold only the "Meta description" label itself (up to, but not including,
# the colon that follows it).
$labelLength = "Meta description".Length
$metaStart = $metaPara.Range.Start
$labelRange = $d.Range($metaStart, $metaStart + $labelLength)
$labelRange.Bold = 1

Write-Host "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
